$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column AA (col 27), shifting old AA/AB -> AC/AD
$ws.Range("AA:AB").Insert()

# New header cells for the inserted columns (order controls sharedStrings index)
$ws.Range("AB1").Value = "Tracking"
$ws.Range("AA1").Value = "ขนส่ง"
